$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename columns
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case small connector words (de/del/la/las/el/los/y) in place names
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("B19").Value = "Amatenango De La Frontera"
$ws.Range("B25").Value = "Chiapa De Corzo"
$ws.Range("B27").Value = "Comitán De Domínguez"
$ws.Range("B50").Value = "Salto De Agua"
$ws.Range("B51").Value = "San Cristóbal De Las Casas"
$ws.Range("B74").Value = "Hidalgo Del Parral"
$ws.Range("B77").Value = "San Francisco De Borja"
$ws.Range("B90").Value = "San Juan De Sabinas"
$ws.Range("A96").Value = "Ciudad De México"
$ws.Range("B112").Value = "Nombre De Dios"
$ws.Range("A115").Value = "Estado De México"
$ws.Range("B115").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B116").Value = "Almoloya De Alquisiras"
$ws.Range("B117").Value = "Almoloya De Juárez"
$ws.Range("B122").Value = "Ecatepec De Morelos"
$ws.Range("B124").Value = "Ixtapan De La Sal"
$ws.Range("B129").Value = "Naucalpan De Juárez"
$ws.Range("B132").Value = "San Felipe Del Progreso"
$ws.Range("B133").Value = "Soyaniquilpan De Juárez"
$ws.Range("B140").Value = "Tenango Del Valle"
$ws.Range("B143").Value = "Tlalnepantla De Baz"
$ws.Range("B147").Value = "Villa De Allende"
$ws.Range("B153").Value = "San Miguel De Allende"
$ws.Range("B154").Value = "Apaseo El Alto"
$ws.Range("B155").Value = "Apaseo El Grande"
$ws.Range("B163").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B166").Value = "Jaral Del Progreso"
$ws.Range("B172").Value = "Purísima Del Rincón"
$ws.Range("B175").Value = "San Diego De La Unión"
$ws.Range("B177").Value = "San Francisco Del Rincón"
$ws.Range("B179").Value = "San Luis De La Paz"
$ws.Range("B180").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B183").Value = "Valle De Santiago"
$ws.Range("B187").Value = "Acapulco De Juárez"
$ws.Range("B188").Value = "Ajuchitlán Del Progreso"
$ws.Range("B191").Value = "Atoyac De Álvarez"
$ws.Range("B192").Value = "Ayutla De Los Libres"
$ws.Range("B194").Value = "Buenavista De Cuéllar"
$ws.Range("B195").Value = "Chilapa De Álvarez"
$ws.Range("B196").Value = "Chilpancingo De Los Bravo"
$ws.Range("B199").Value = "Coyuca De Benítez"
$ws.Range("B200").Value = "Coyuca De Catalán"
$ws.Range("B203").Value = "Huitzuco De Los Figueroa"
$ws.Range("B204").Value = "Iguala De La Independencia"
$ws.Range("B205").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B211").Value = "Taxco De Alarcón"
$ws.Range("B213").Value = "Técpan De Galeana"
$ws.Range("B215").Value = "Tixtla De Guerrero"
$ws.Range("B226").Value = "Atotonilco De Tula"
$ws.Range("B230").Value = "Cuautepec De Hinojosa"
$ws.Range("B234").Value = "Huejutla De Reyes"
$ws.Range("B237").Value = "Jacala De Ledezma"
$ws.Range("B240").Value = "Molango De Escamilla"
$ws.Range("B241").Value = "Pachuca De Soto"
$ws.Range("B247").Value = "Tenango De Doria"
$ws.Range("B249").Value = "Tepehuacán De Guerrero"
$ws.Range("B250").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B254").Value = "Tula De Allende"
$ws.Range("B256").Value = "Zacualtipán De Ángeles"
$ws.Range("B261").Value = "Autlán De Navarro"
$ws.Range("B270").Value = "Jilotlán De Los Dolores"
$ws.Range("B271").Value = "Lagos De Moreno"
$ws.Range("B275").Value = "San Juanito De Escobedo"
$ws.Range("B277").Value = "Tamazula De Gordiano"
$ws.Range("B278").Value = "Tepatitlán De Morelos"
$ws.Range("B280").Value = "Tizapán El Alto"
$ws.Range("B283").Value = "Zacoalco De Torres"
$ws.Range("B328").Value = "Tetela Del Volcán"
$ws.Range("B342").Value = "Mier Y Noriega"
$ws.Range("B345").Value = "San Nicolás De Los Garza"
$ws.Range("B348").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B349").Value = "Coicoyán De Las Flores"
$ws.Range("B351").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B352").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B353").Value = "Ixtlán De Juárez"
$ws.Range("B354").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B357").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B358").Value = "Oaxaca De Juárez"
$ws.Range("B359").Value = "Putla Villa De Guerrero"
$ws.Range("B360").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B370").Value = "San Miguel Del Puerto"
$ws.Range("B381").Value = "Santa Inés De Zaragoza"
$ws.Range("B382").Value = "Santa Inés Del Monte"
$ws.Range("B402").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B403").Value = "Tlacolula De Matamoros"
$ws.Range("B404").Value = "Villa De Etla"
$ws.Range("B405").Value = "Villa De Tututepec"
$ws.Range("B406").Value = "Villa Sola De Vega"
$ws.Range("B407").Value = "Zimatlán De Álvarez"
$ws.Range("B417").Value = "Cuayuca De Andrade"
$ws.Range("B424").Value = "Izúcar De Matamoros"
$ws.Range("B426").Value = "Los Reyes De Juárez"
$ws.Range("B428").Value = "Palmar De Bravo"
$ws.Range("B441").Value = "Amealco De Bonfil"
$ws.Range("B443").Value = "Cadereyta De Montes"
$ws.Range("B447").Value = "Jalpan De Serra"
$ws.Range("B448").Value = "Landa De Matamoros"
$ws.Range("B449").Value = "Pinal De Amoles"
$ws.Range("B452").Value = "San Juan Del Río"
$ws.Range("B458").Value = "Armadillo De Los Infante"
$ws.Range("B459").Value = "Axtla De Terrazas"
$ws.Range("B462").Value = "Cerro De San Pedro"
$ws.Range("B463").Value = "Ciudad Del Maíz"
$ws.Range("B470").Value = "Mexquitic De Carmona"
$ws.Range("B474").Value = "San Ciro De Acosta"
$ws.Range("B478").Value = "Santa María Del Río"
$ws.Range("B485").Value = "Villa De Arriaga"
$ws.Range("B486").Value = "Villa De Ramos"
$ws.Range("B487").Value = "Villa De Reyes"
$ws.Range("B529").Value = "Soto La Marina"
$ws.Range("B540").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B550").Value = "Amatlán De Los Reyes"
$ws.Range("B555").Value = "Boca Del Río"
$ws.Range("B565").Value = "Cosamaloapan De Carpio"
$ws.Range("B575").Value = "Hueyapan De Ocampo"
$ws.Range("B586").Value = "Juchique De Ferrer"
$ws.Range("B589").Value = "Martínez De La Torre"
$ws.Range("B599").Value = "Paso De Ovejas"
$ws.Range("B602").Value = "Poza Rica De Hidalgo"
$ws.Range("B608").Value = "Sayula De Alemán"
$ws.Range("B610").Value = "Soledad De Doblado"
$ws.Range("B627").Value = "Vega De Alatorre"
$ws.Range("B643").Value = "Noria De Ángeles"
$ws.Range("B650").Value = "Villa De Cos"

# Floating point precision correction for percentage column
$ws.Range("D10").Value = 0.0009250693802035152
$ws.Range("D23").Value = 0.0009250693802035152
$ws.Range("D34").Value = 0.0009250693802035152
$ws.Range("D44").Value = 0.0009250693802035152
$ws.Range("D47").Value = 0.0009250693802035152
$ws.Range("D48").Value = 0.0009250693802035152
$ws.Range("D62").Value = 0.0009250693802035152
$ws.Range("D65").Value = 0.0009250693802035152
$ws.Range("D72").Value = 0.0009250693802035152
$ws.Range("D74").Value = 0.0009250693802035152
$ws.Range("D76").Value = 0.0009250693802035152
$ws.Range("D88").Value = 0.0009250693802035152
$ws.Range("D95").Value = 0.0009250693802035152
$ws.Range("D106").Value = 0.0009250693802035152
$ws.Range("D117").Value = 0.0009250693802035152
$ws.Range("D121").Value = 0.0009250693802035152
$ws.Range("D122").Value = 0.0009250693802035152
$ws.Range("D126").Value = 0.0009250693802035152
$ws.Range("D137").Value = 0.0009250693802035152
$ws.Range("D141").Value = 0.0009250693802035152
$ws.Range("D148").Value = 0.0009250693802035152
$ws.Range("D149").Value = 0.0009250693802035152
$ws.Range("D151").Value = 0.0009250693802035152
$ws.Range("D160").Value = 0.0009250693802035152
$ws.Range("D162").Value = 0.0009250693802035152
$ws.Range("D166").Value = 0.0009250693802035152
$ws.Range("D171").Value = 0.0009250693802035152
$ws.Range("D174").Value = 0.0009250693802035152
$ws.Range("D178").Value = 0.0009250693802035152
$ws.Range("D192").Value = 0.0009250693802035152
$ws.Range("D195").Value = 0.0009250693802035152
$ws.Range("D200").Value = 0.0009250693802035152
$ws.Range("D201").Value = 0.0009250693802035152
$ws.Range("D204").Value = 0.0009250693802035152
$ws.Range("D205").Value = 0.0009250693802035152
$ws.Range("D212").Value = 0.0009250693802035152
$ws.Range("D213").Value = 0.0009250693802035152
$ws.Range("D218").Value = 0.0009250693802035152
$ws.Range("D230").Value = 0.0009250693802035152
$ws.Range("D234").Value = 0.0009250693802035152
$ws.Range("D235").Value = 0.0009250693802035152
$ws.Range("D237").Value = 0.0009250693802035152
$ws.Range("D242").Value = 0.0009250693802035152
$ws.Range("D256").Value = 0.0009250693802035152
$ws.Range("D259").Value = 0.0009250693802035152
$ws.Range("D272").Value = 0.0009250693802035152
$ws.Range("D282").Value = 0.0009250693802035152
$ws.Range("D284").Value = 0.0009250693802035152
$ws.Range("D285").Value = 0.0009250693802035152
$ws.Range("D290").Value = 0.0009250693802035152
$ws.Range("D297").Value = 0.0009250693802035152
$ws.Range("D299").Value = 0.0009250693802035152
$ws.Range("D308").Value = 0.0009250693802035152
$ws.Range("D311").Value = 0.0009250693802035152
$ws.Range("D314").Value = 0.0009250693802035152
$ws.Range("D315").Value = 0.0009250693802035152
$ws.Range("D333").Value = 0.0009250693802035152
$ws.Range("D334").Value = 0.0009250693802035152
$ws.Range("D338").Value = 0.0009250693802035152
$ws.Range("D340").Value = 0.0009250693802035152
$ws.Range("D342").Value = 0.0009250693802035152
$ws.Range("D345").Value = 0.0009250693802035152
$ws.Range("D349").Value = 0.0009250693802035152
$ws.Range("D351").Value = 0.0009250693802035152
$ws.Range("D352").Value = 0.0009250693802035152
$ws.Range("D355").Value = 0.0009250693802035152
$ws.Range("D356").Value = 0.0009250693802035152
$ws.Range("D358").Value = 0.0009250693802035152
$ws.Range("D359").Value = 0.0009250693802035152
$ws.Range("D360").Value = 0.0009250693802035152
$ws.Range("D361").Value = 0.0009250693802035152
$ws.Range("D363").Value = 0.0009250693802035152
$ws.Range("D366").Value = 0.0009250693802035152
$ws.Range("D369").Value = 0.0009250693802035152
$ws.Range("D376").Value = 0.0009250693802035152
$ws.Range("D380").Value = 0.0009250693802035152
$ws.Range("D384").Value = 0.0009250693802035152
$ws.Range("D386").Value = 0.0009250693802035152
$ws.Range("D390").Value = 0.0009250693802035152
$ws.Range("D395").Value = 0.0009250693802035152
$ws.Range("D401").Value = 0.0009250693802035152
$ws.Range("D405").Value = 0.0009250693802035152
$ws.Range("D407").Value = 0.0009250693802035152
$ws.Range("D415").Value = 0.0009250693802035152
$ws.Range("D416").Value = 0.0009250693802035152
$ws.Range("D418").Value = 0.0009250693802035152
$ws.Range("D422").Value = 0.0009250693802035152
$ws.Range("D428").Value = 0.0009250693802035152
$ws.Range("D433").Value = 0.0009250693802035152
$ws.Range("D446").Value = 0.0009250693802035152
$ws.Range("D449").Value = 0.0009250693802035152
$ws.Range("D456").Value = 0.0009250693802035152
$ws.Range("D463").Value = 0.0009250693802035152
$ws.Range("D467").Value = 0.0009250693802035152
$ws.Range("D469").Value = 0.0009250693802035152
$ws.Range("D470").Value = 0.0009250693802035152
$ws.Range("D487").Value = 0.0009250693802035152
$ws.Range("D496").Value = 0.0009250693802035152
$ws.Range("D505").Value = 0.0009250693802035152
$ws.Range("D506").Value = 0.0009250693802035152
$ws.Range("D516").Value = 0.0009250693802035152
$ws.Range("D530").Value = 0.0009250693802035152
$ws.Range("D536").Value = 0.0009250693802035152
$ws.Range("D538").Value = 0.0009250693802035152
$ws.Range("D541").Value = 0.0009250693802035152
$ws.Range("D546").Value = 0.0009250693802035152
$ws.Range("D587").Value = 0.0009250693802035152
$ws.Range("D611").Value = 0.0009250693802035152
$ws.Range("D621").Value = 0.0009250693802035152
$ws.Range("D622").Value = 0.0009250693802035152
$ws.Range("D626").Value = 0.0009250693802035152
$ws.Range("D630").Value = 0.0009250693802035152
$ws.Range("D650").Value = 0.0009250693802035152

# Remove trailing footer/metadata rows (657-661)
$ws.Rows.Item(657).Resize(5).Delete()

Write-Host "edit complete"